# tracks split identifiers column
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Identifiers" column header to "Internal House Name"
$ws.Range("O1").Value = "Internal House Name"

# New column for the split-off identifier: Tag/Band
$ws.Range("S1").Value = "Tag /Band"

# Match the selection/viewport shown in the diff after the edit
$ws.Range("S1").Select()
